$wb = $excel.ActiveWorkbook

$sTask = $wb.Worksheets.Item("TASK SUMMARY SHEET")
$sSummary = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

# --- Chronological order matters: Excel appends new unique strings to the
# shared-string table in the order they are first typed. Replicate that
# order so the resulting file matches what the author actually produced.

# 1) Summary sheet: name of author
$sSummary.Range("D1").Value = "Jesse"

# 2) Task summary sheet: first row (stage/task already in progress before this update)
$sTask.Range("A3").Value = "Proj Design Specifics"

# 3) Task summary sheet header: author's full name + week number
$sTask.Range("C1").Value = "Jesse Hare"
$sTask.Range("E1").Value = 2

# 4) Task summary sheet: "Project Build" stage rows
$sTask.Range("A4").Value = "Project Build"
$sTask.Range("A5").Value = "Project Build"
$sTask.Range("A6").Value = "Project Build"
$sTask.Range("A7").Value = "Project Build"

# 5) Task descriptions for the Project Build stage
$sTask.Range("B4").Value = "Build front end GUI first iteration "
$sTask.Range("B5").Value = "Write code for search algorithm "
$sTask.Range("B6").Value = "Write code for sorting CSV imported elements based on user input "

# 6) Task description for the first (Proj Design Specifics) row
$sTask.Range("B3").Value = "Finalise framework choice and libraries needed"

# 7) Summary sheet stage name
$sSummary.Range("A4").Value = "Project Design and Specifics"

# 8) Last task row description
$sTask.Range("B7").Value = "Iteration Review"

# --- Numeric data for Task Summary Sheet ---
$sTask.Range("C3").Value = 3
$sTask.Range("D3").Value = 2
$sTask.Range("E3").Value = 0

$sTask.Range("C4").Value = 20
$sTask.Range("D4").Value = 20
$sTask.Range("E4").Value = 0

$sTask.Range("C5").Value = 20
$sTask.Range("D5").Value = 30
$sTask.Range("E5").Value = 0

$sTask.Range("C6").Value = 20
$sTask.Range("D6").Value = 30
$sTask.Range("E6").Value = 0

$sTask.Range("C7").Value = 2
$sTask.Range("D7").Value = 1
$sTask.Range("E7").Value = 0

# --- Numeric data for Activity Log Summary Sheet ---
$sSummary.Range("B4").Value = 2
$sSummary.Range("A5").Value = "Project Build"
$sSummary.Range("B5").Value = 81

# --- Active sheet / selection updates ---
# Leave the Activity Log Summary Sheet with B6 selected (last place the author
# clicked there before switching back to the Task Summary Sheet).
[void]$sSummary.Activate()
[void]$sSummary.Range("B6").Select()

[void]$sTask.Activate()
[void]$sTask.Range("E7").Select()
